# Update column F (dSF) values on Sheet1 to reflect re-pulled data / mean calc.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("F3").Value = 2
$ws.Range("F7").Value = -3
$ws.Range("F8").Value = -6
$ws.Range("F16").Value = 10
$ws.Range("F17").Value = 5
$ws.Range("F21").Value = -2
$ws.Range("F23").Value = 3
$ws.Range("F27").Value = 1
$ws.Range("F33").Value = 0
$ws.Range("F41").Value = 1
$ws.Range("F43").Value = 1
$ws.Range("F44").Value = 4
$ws.Range("F46").Value = 4
$ws.Range("F49").Value = -1
$ws.Range("F50").Value = -13
$ws.Range("F51").Value = -3
$ws.Range("F52").Value = 1
$ws.Range("F53").Value = -6
$ws.Range("F55").Value = 10
$ws.Range("F61").Value = -4
$ws.Range("F65").Value = -1
